# Update Active_Outages.xlsx - 6/18/2025, 5:05:05 PM
#
# Applies the following per-cell updates (Elapsed Duration(Hrs) refresh plus
# two site-status corrections) across the R1..R6 sheets.

$wb = $excel.ActiveWorkbook

# --- R1 ---
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3930:19:13"
$ws.Range("G3").Value = "69:51:51"
$ws.Range("G4").Value = "92:51:51"
$ws.Range("D6").Value = "HAJ0155"
$ws.Range("J6").Value = "Good"

# --- R2 ---
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12111:42:52"
$ws.Range("G3").Value = "3241:26:21"
$ws.Range("G4").Value = "479:37:55"
$ws.Range("D5").Value = "JED0190"
$ws.Range("J5").Value = "Good"

# --- R4 ---
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2957:32:41"
$ws.Range("G3").Value = "184:44:56"
$ws.Range("G4").Value = "72:57:21"
$ws.Range("G5").Value = "70:34:54"

# --- R5 ---
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "431:31:40"

# --- R6 ---
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "72:03:58"
